$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.211.15"
$ws.Range("E2").Value = "  +5.65%  "
$ws.Range("D3").Value = "2.463.56"
$ws.Range("E3").Value = "  +6.86%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "566.64"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +4.27%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.47"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +10.99%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  +2.82%  "
$ws.Range("D9").Value = "2.462.68"
$ws.Range("E9").Value = "  +6.85%  "
$ws.Range("E10").Value = "  +5.81%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.72"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +3.18%  "
$ws.Range("E12").Value = "  +1.29%  "
$ws.Range("E13").Value = "  +5.94%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.52"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +13.85%  "
$ws.Range("D15").Value = "2.903.32"
$ws.Range("E15").Value = "  +6.75%  "
$ws.Range("D16").Value = "63.119.43"
$ws.Range("E16").Value = "  +5.36%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000143"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +7.07%  "
$ws.Range("D18").Value = "2.463.68"
$ws.Range("E18").Value = "  +6.44%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.23"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +6.88%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "341.85"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +9.66%  "
$ws.Range("E21").Value = "  +6.22%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.81"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +4.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.998"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.21%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.62"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +2.72%  "
$ws.Range("E25").Value = "  +2.80%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("E27").Value = "  +11.14%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.15"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +4.94%  "
$ws.Range("E29").Value = "  +8.89%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.85"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +15.43%  "
$ws.Range("D31").Value = "0.0₃0814"
$ws.Range("E31").Value = "  +12.93%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.84"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +7.44%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "175.04"
$ws.Range("D33").ClearFormats()
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.52"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +12.66%  "
$ws.Range("E35").Value = "  +4.94%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.93"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +5.62%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "369.84"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +17.65%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.46"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +10.33%  "
$ws.Range("E39").Value = "  -0.01%  "
$ws.Range("E40").Value = "  -0.11%  "
$ws.Range("E41").Value = "  +13.27%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "40.42"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +6.58%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "152.07"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +11.30%  "
$ws.Range("E44").Value = "  +7.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "20.57"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +8.68%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.599"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +6.73%  "
$ws.Range("E47").Value = "  +2.76%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0520"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +5.36%  "
$ws.Range("D49").Value = "0.0₆0240"
$ws.Range("E49").Value = "  +11.33%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0226"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +5.58%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "18.01"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +7.34%  "
